$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells we rewrite keep text formatting like the source data
$ws.Range("D2:D9").NumberFormat = "@"
$ws.Range("D11:D13").NumberFormat = "@"
$ws.Range("D15:D42").NumberFormat = "@"
$ws.Range("D44:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "23.735.89"
$ws.Range("E2").Value = "  +1.81%  "

# Row 3
$ws.Range("D3").Value = "1.653.16"
$ws.Range("E3").Value = "  +1.72%  "

# Row 4
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "0.9996"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").Value = "303.73"
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("D7").Value = "0.3826"
$ws.Range("E7").Value = "  +2.16%  "

# Row 8
$ws.Range("D8").Value = "51.32"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "0.3609"
$ws.Range("E9").Value = "  -0.32%  "

# Row 10
$ws.Range("E10").Value = "  +2.12%  "

# Row 11
$ws.Range("D11").Value = "0.08235"
$ws.Range("E11").Value = "  +1.14%  "

# Row 12
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.24%  "

# Row 13
$ws.Range("D13").Value = "22.64"
$ws.Range("E13").Value = "  +1.84%  "

# Row 14
$ws.Range("E14").Value = "  +1.18%  "

# Row 15
$ws.Range("D15").Value = "7.412"
$ws.Range("E15").Value = "  +1.60%  "

# Row 16
$ws.Range("D16").Value = "0.00001235"
$ws.Range("E16").Value = "  -0.18%  "

# Row 17
$ws.Range("D17").Value = "1.648.88"
$ws.Range("E17").Value = "  +1.93%  "

# Row 18
$ws.Range("D18").Value = "97.44"
$ws.Range("E18").Value = "  +3.82%  "

# Row 19
$ws.Range("D19").Value = "0.06977"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("D20").Value = "6.776"
$ws.Range("E20").Value = "  +3.44%  "

# Row 21
$ws.Range("D21").Value = "17.69"
$ws.Range("E21").Value = "  +1.14%  "

# Row 22
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").Value = "12.63"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24
$ws.Range("D24").Value = "23.725.25"
$ws.Range("E24").Value = "  +1.72%  "

# Row 25
$ws.Range("D25").Value = "2.527"
$ws.Range("E25").Value = "  +2.29%  "

# Row 26
$ws.Range("D26").Value = "3.087"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("D27").Value = "21.32"
$ws.Range("E27").Value = "  +0.83%  "

# Row 28
$ws.Range("D28").Value = "151.60"
$ws.Range("E28").Value = "  +0.97%  "

# Row 29
$ws.Range("D29").Value = "5.276"
$ws.Range("E29").Value = "  +0.68%  "

# Row 30
$ws.Range("D30").Value = "134.97"
$ws.Range("E30").Value = "  +1.83%  "

# Row 31
$ws.Range("D31").Value = "1.835.09"
$ws.Range("E31").Value = "  +1.84%  "

# Row 32
$ws.Range("D32").Value = "6.868"
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("D33").Value = "1.093"
$ws.Range("E33").Value = "  +5.33%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.109"
$ws.Range("E34").Value = "  +2.19%  "

# Row 35
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "11.83"
$ws.Range("E35").Value = "  +9.60%  "

# Row 36
$ws.Range("D36").Value = "0.02845"
$ws.Range("E36").Value = "  +3.32%  "

# Row 37
$ws.Range("D37").Value = "0.2518"
$ws.Range("E37").Value = "  +1.08%  "

# Row 38
$ws.Range("D38").Value = "0.08831"
$ws.Range("E38").Value = "  +0.60%  "

# Row 39
$ws.Range("D39").Value = "6.101"
$ws.Range("E39").Value = "  +2.10%  "

# Row 40
$ws.Range("D40").Value = "0.07049"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
$ws.Range("D41").Value = "12.88"
$ws.Range("E41").Value = "  +6.79%  "

# Row 42
$ws.Range("D42").Value = "0.7078"
$ws.Range("E42").Value = "  +1.51%  "

# Row 43
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").Value = "16.02"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45
$ws.Range("D45").Value = "0.6551"
$ws.Range("E45").Value = "  +1.22%  "

# Row 46
$ws.Range("E46").Value = "  +3.45%  "

# Row 47
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
$ws.Range("D48").Value = "3.972"
$ws.Range("E48").Value = "  +0.27%  "

# Row 49
$ws.Range("D49").Value = "0.07989"
$ws.Range("E49").Value = "  +0.33%  "

# Row 50
$ws.Range("D50").Value = "128.50"
$ws.Range("E50").Value = "  +2.49%  "

# Row 51
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  +1.09%  "
